# Adds the ability to edit the variations from the editor.
# Appends two new parameter rows (Variants tab) to the game_params sheet,
# and updates the sheet view (frozen-pane scroll position / selection)
# to match where the editor was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("game_params")

# -- New row 87: "variants" (TextDict of Named Variants) -------------------
$ws.Range("A87").Value = "Variants"
$ws.Range("B87").Value = "variants"
$ws.Range("C87").Value = "Named Variants"
$ws.Range("D87").Value = "_"
$ws.Range("E87").Value = 250
$ws.Range("F87").Value = "TextDict"
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0

# -- New row 88: "vari_params" (TextDict of Vari Params) -------------------
$ws.Range("A88").Value = "Variants"
$ws.Range("B88").Value = "vari_params"
$ws.Range("C88").Value = "Vari Params"
$ws.Range("D88").Value = "_"
$ws.Range("E88").Value = 251
$ws.Range("F88").Value = "TextDict"
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 2

# -- Sheet view: keep the header row frozen, scroll/select near the new rows
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 73
$ws.Range("A87").Select()
